$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.964.89'
$ws.Range('E2').Value = '  +0.41%  '
$ws.Range('D3').Value = '2.582.14'
$ws.Range('E3').Value = '  +2.54%  '
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').Value = '''302.88'
$ws.Range('E5').Value = '  +2.15%  '
$ws.Range('D6').Value = '''97.38'
$ws.Range('E6').Value = '  +4.54%  '
$ws.Range('D7').Value = '''0.576'
$ws.Range('E7').Value = '  +1.33%  '
$ws.Range('E8').Value = '  +0.14%  '
$ws.Range('D9').Value = '''0.552'
$ws.Range('E9').Value = '  +0.91%  '
$ws.Range('D10').Value = '''36.93'
$ws.Range('E10').Value = '  +1.86%  '
$ws.Range('D11').Value = '''0.0810'
$ws.Range('E11').Value = '  +1.46%  '
$ws.Range('D12').Value = '''7.75'
$ws.Range('E12').Value = '  +1.61%  '
$ws.Range('E13').Value = '  +7.53%  '
$ws.Range('D14').Value = '2.596.02'
$ws.Range('E14').Value = '  +3.25%  '
$ws.Range('D15').Value = '''0.890'
$ws.Range('E15').Value = '  +2.86%  '
$ws.Range('D16').Value = '''14.41'
$ws.Range('E16').Value = '  +3.10%  '
$ws.Range('D17').Value = '42.989.78'
$ws.Range('E17').Value = '  +0.64%  '
$ws.Range('D18').Value = '''12.97'
$ws.Range('E18').Value = '  +5.69%  '
$ws.Range('D19').Value = '0.0₃0993'
$ws.Range('E19').Value = '  +3.31%  '
$ws.Range('D20').Value = '''6.68'
$ws.Range('E20').Value = '  +2.64%  '
$ws.Range('D21').Value = '''72.15'
$ws.Range('E21').Value = '  -0.31%  '
$ws.Range('D22').Value = '''255.34'
$ws.Range('E22').Value = '  -1.29%  '
$ws.Range('D23').Value = '''2.97'
$ws.Range('E23').Value = '  +3.17%  '
$ws.Range('D24').Value = '''2.14'
$ws.Range('E24').Value = '  +0.63%  '
$ws.Range('D25').Value = '''28.68'
$ws.Range('E25').Value = '  -1.15%  '
$ws.Range('E26').Value = '  -0.09%  '
$ws.Range('D27').Value = '''10.24'
$ws.Range('E27').Value = '  +3.31%  '
$ws.Range('D28').Value = '''38.81'
$ws.Range('E28').Value = '  +5.34%  '
$ws.Range('D30').Value = '''6.07'
$ws.Range('E30').Value = '  +1.36%  '
$ws.Range('D31').Value = '''155.68'
$ws.Range('E31').Value = '  +3.32%  '
$ws.Range('D32').Value = '''2.19'
$ws.Range('E32').Value = '  +0.12%  '
$ws.Range('D33').Value = '''2.75'
$ws.Range('E33').Value = '  -0.35%  '
$ws.Range('E34').Value = '  +2.56%  '
$ws.Range('E35').Value = '  -2.31%  '
$ws.Range('D36').Value = '''18.44'
$ws.Range('E36').Value = '  +13.47%  '
$ws.Range('E37').Value = '  +0.92%  '
$ws.Range('E38').Value = '  +1.63%  '
$ws.Range('D39').Value = '''23.39'
$ws.Range('E39').Value = '  -1.89%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D40').Value = '''3.91'
$ws.Range('E40').Value = '  +2.47%  '
$ws.Range('B41').Value = 'NEARProtocol'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D41').Value = '''3.43'
$ws.Range('E41').Value = '  +0.43%  '
$ws.Range('E42').Value = '  +1.10%  '
$ws.Range('E43').Value = '  +29.09%  '
$ws.Range('D44').Value = '2.066.17'
$ws.Range('E44').Value = '  +2.73%  '
$ws.Range('E45').Value = '  +0.38%  '
$ws.Range('D46').Value = '''9.26'
$ws.Range('E46').Value = '  +4.51%  '
$ws.Range('D47').Value = '''85.61'
$ws.Range('E47').Value = '  -0.13%  '
$ws.Range('D48').Value = '''77.03'
$ws.Range('E48').Value = '  +13.21%  '
$ws.Range('B49').Value = 'RocketPoolETH'
$ws.Range('C49').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D49').Value = '2.829.84'
$ws.Range('E49').Value = '  +2.91%  '
$ws.Range('B50').Value = 'Aave'
$ws.Range('C50').Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range('D50').Value = '''106.44'
$ws.Range('E50').Value = '  +3.61%  '
$ws.Range('B51').Value = 'Stacks'
$ws.Range('C51').Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range('D51').Value = '''1.69'
$ws.Range('E51').Value = '  +2.89%  '
